$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.656.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.862.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3911"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.33"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07961"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.859.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.949"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.216"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06733"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001042"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.643.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  -1.08%  "
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.309"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.082.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.131"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.427"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9776"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09408"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.628"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.306"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.332"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02233"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06019"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.346"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.193"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5938"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1869"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5608"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06717"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.052"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.59%  "
